$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 190.16667
$ws.Range("J4").Value = 195
$ws.Range("L4").Value = 195
$ws.Range("N4").Value = -423

$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 3000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2931

$ws.Range("H100").Value = 1883.8422
$ws.Range("I100").Value = 1645.2307
$ws.Range("K100").Value = 1645.2307
$ws.Range("M100").Value = -1104.2307

$ws.Range("H113").Value = 4330.4443
$ws.Range("I113").Value = 4163.5
$ws.Range("J113").Value = 4664.3335
$ws.Range("K113").Value = 4163.5
$ws.Range("L113").Value = 4664.3335
$ws.Range("M113").Value = -909.5
$ws.Range("N113").Value = -11172.3335

$ws.Range("H118").Value = 1015.4
$ws.Range("I118").Value = 519.25
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 1557.75
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 99.25
$ws.Range("N118").Value = -12314

$ws.Range("H129").Value = 1776.4166
$ws.Range("I129").Value = 1058.1666
$ws.Range("K129").Value = 3174.4998
$ws.Range("M129").Value = 1825.5002

$ws.Range("H135").Value = 922.75
$ws.Range("I135").Value = 790.2778
$ws.Range("K135").Value = 7112.500199999999
$ws.Range("M135").Value = -4577.500199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 12130.5
$ws.Range("J44").Value = 12130.5
$ws.Range("L44").Value = 12130.5
$ws.Range("N44").Value = -13106.5

$ws.Range("H61").Value = 1269.3889
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H122").Value = 2484.3076
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 11250
$ws.Range("N122").Value = -16150

$ws.Range("H136").Value = 1269.3889
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1173.1111
$ws.Range("I99").Value = 1197.7142
$ws.Range("J99").Value = 1087
$ws.Range("K99").Value = 1197.7142
$ws.Range("L99").Value = 1087
$ws.Range("M99").Value = 300.2858000000001
$ws.Range("N99").Value = -4083

$ws.Range("H107").Value = 4847.1875
$ws.Range("I107").Value = 875.25
$ws.Range("K107").Value = 875.25
$ws.Range("M107").Value = 1044.75

$ws.Range("H134").Value = 1516.1818
$ws.Range("I134").Value = 1317
$ws.Range("K134").Value = 3951
$ws.Range("M134").Value = -1416

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 54.363636
$ws.Range("K7").Value = 54.363636
$ws.Range("M7").Value = 58.636364

$ws.Range("H122").Value = 941.63635
$ws.Range("I122").Value = 885.8
$ws.Range("K122").Value = 2657.4
$ws.Range("M122").Value = -207.3999999999996

$ws.Range("H132").Value = 3785.7896
$ws.Range("I132").Value = 3587.7334
$ws.Range("J132").Value = 4528.5
$ws.Range("K132").Value = 10763.2002
$ws.Range("L132").Value = 13585.5
$ws.Range("M132").Value = -8233.200199999999
$ws.Range("N132").Value = -18645.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2107.0715
$ws.Range("I57").Value = 1999.909
$ws.Range("J57").Value = 2500
$ws.Range("K57").Value = 5999.727000000001
$ws.Range("L57").Value = 7500
$ws.Range("M57").Value = -5440.727000000001
$ws.Range("N57").Value = -8618

$ws.Range("H131").Value = 2381.889
$ws.Range("I131").Value = 1509.8
$ws.Range("K131").Value = 4529.4
$ws.Range("M131").Value = 510.6000000000004

$ws.Range("H138").Value = 4790.9
$ws.Range("J138").Value = 7141.6665
$ws.Range("L138").Value = 21424.9995
$ws.Range("N138").Value = -31704.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H113").Value = 7316.0713
$ws.Range("I113").Value = 4080.4
$ws.Range("J113").Value = 9113.666999999999
$ws.Range("K113").Value = 4080.4
$ws.Range("L113").Value = 9113.666999999999
$ws.Range("M113").Value = -1910.4
$ws.Range("N113").Value = -13453.667

$ws.Range("H122").Value = 1868.875
$ws.Range("I122").Value = 1909.1666
$ws.Range("J122").Value = 1748
$ws.Range("K122").Value = 5727.4998
$ws.Range("L122").Value = 5244
$ws.Range("M122").Value = -3277.4998
$ws.Range("N122").Value = -10144

$ws.Range("H126").Value = 3205.375
$ws.Range("I126").Value = 2949
$ws.Range("K126").Value = 8847
$ws.Range("M126").Value = -6377

$ws.Range("H132").Value = 49071.637
$ws.Range("I132").Value = 68968.53
$ws.Range("K132").Value = 206905.59
$ws.Range("M132").Value = -204375.59

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 839.63635
$ws.Range("J22").Value = 838
$ws.Range("L22").Value = 838
$ws.Range("N22").Value = -1428

$ws.Range("H27").Value = 839.63635
$ws.Range("J27").Value = 838
$ws.Range("L27").Value = 838
$ws.Range("N27").Value = -1052

$ws.Range("H46").Value = 3330.4119
$ws.Range("I46").Value = 2077.8
$ws.Range("J46").Value = 3852.3333
$ws.Range("K46").Value = 2077.8
$ws.Range("L46").Value = 3852.3333
$ws.Range("M46").Value = -1889.8
$ws.Range("N46").Value = -4228.3333

$ws.Range("H55").Value = 1593.2222
$ws.Range("I55").Value = 1939.8334
$ws.Range("K55").Value = 1939.8334
$ws.Range("M55").Value = -1766.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 797.93335
$ws.Range("J96").Value = 642.5
$ws.Range("L96").Value = 642.5
$ws.Range("N96").Value = -3388.5

$ws.Range("H107").Value = 735.8
$ws.Range("I107").Value = 297
$ws.Range("K107").Value = 891
$ws.Range("M107").Value = 1029

$ws.Range("H122").Value = 3727.1
$ws.Range("I122").Value = 1570.1666
$ws.Range("K122").Value = 4710.4998
$ws.Range("M122").Value = -2260.4998

$ws.Range("H126").Value = 3309.4167
$ws.Range("I126").Value = 1819.7646
$ws.Range("J126").Value = 6927.143
$ws.Range("K126").Value = 5459.293799999999
$ws.Range("L126").Value = 20781.429
$ws.Range("M126").Value = -2989.293799999999
$ws.Range("N126").Value = -25721.429

$ws.Range("H132").Value = 1974
$ws.Range("I132").Value = 2044.3636
$ws.Range("K132").Value = 6133.0908
$ws.Range("M132").Value = -3603.0908
